# Auto-generated edit script applying scheduled-runner profit updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716
$ws.Range("H94").Value = 540
$ws.Range("I94").Value = 540
$ws.Range("K94").Value = 540
$ws.Range("M94").Value = -89
$ws.Range("H116").Value = 3849.8
$ws.Range("J116").Value = 3899.3333
$ws.Range("L116").Value = 3899.3333
$ws.Range("N116").Value = -10783.3333
$ws.Range("H129").Value = 1705.5
$ws.Range("I129").Value = 1447.25
$ws.Range("J129").Value = 2222
$ws.Range("K129").Value = 4341.75
$ws.Range("L129").Value = 6666
$ws.Range("M129").Value = 658.25
$ws.Range("N129").Value = -16666
$ws.Range("H132").Value = 27449.5
$ws.Range("I132").Value = 27449.5
$ws.Range("K132").Value = 82348.5
$ws.Range("M132").Value = -79818.5
$ws.Range("H137").Value = 1646.7
$ws.Range("J137").Value = 2666.6667
$ws.Range("L137").Value = 8000.000100000001
$ws.Range("N137").Value = -13100.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3267.9412
$ws.Range("I2").Value = 999.5
$ws.Range("J2").Value = 5284.3335
$ws.Range("K2").Value = 999.5
$ws.Range("L2").Value = 5284.3335
$ws.Range("M2").Value = -886.5
$ws.Range("N2").Value = -5510.3335
$ws.Range("H46").Value = 6999.1816
$ws.Range("J46").Value = 7232.3335
$ws.Range("L46").Value = 7232.3335
$ws.Range("N46").Value = -7870.3335
$ws.Range("H74").Value = 2935.743
$ws.Range("I74").Value = 2650.0303
$ws.Range("J74").Value = 7650
$ws.Range("K74").Value = 2650.0303
$ws.Range("L74").Value = 7650
$ws.Range("M74").Value = -1776.0303
$ws.Range("N74").Value = -9398
$ws.Range("H77").Value = 2935.743
$ws.Range("I77").Value = 2650.0303
$ws.Range("J77").Value = 7650
$ws.Range("K77").Value = 13250.1515
$ws.Range("L77").Value = 38250
$ws.Range("M77").Value = -8882.1515
$ws.Range("N77").Value = -46986
$ws.Range("H116").Value = 3267.9412
$ws.Range("I116").Value = 999.5
$ws.Range("J116").Value = 5284.3335
$ws.Range("K116").Value = 999.5
$ws.Range("L116").Value = 5284.3335
$ws.Range("M116").Value = 1294.5
$ws.Range("N116").Value = -9872.333500000001
$ws.Range("H132").Value = 3314.5557
$ws.Range("I132").Value = 2453.641
$ws.Range("K132").Value = 7360.923000000001
$ws.Range("M132").Value = -4830.923000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3267.9412
$ws.Range("I3").Value = 999.5
$ws.Range("J3").Value = 5284.3335
$ws.Range("K3").Value = 999.5
$ws.Range("L3").Value = 5284.3335
$ws.Range("M3").Value = -885.5
$ws.Range("N3").Value = -5512.3335
$ws.Range("H59").Value = 199999
$ws.Range("J59").Value = 199999
$ws.Range("L59").Value = 199999
$ws.Range("N59").Value = -201693
$ws.Range("H134").Value = 3374.9211
$ws.Range("I134").Value = 3412.4849
$ws.Range("K134").Value = 10237.4547
$ws.Range("M134").Value = -7702.4547

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4167.7144
$ws.Range("I31").Value = 4270.857
$ws.Range("J31").Value = 4098.952
$ws.Range("K31").Value = 4270.857
$ws.Range("L31").Value = 4098.952
$ws.Range("M31").Value = -3975.857
$ws.Range("N31").Value = -4688.952
$ws.Range("H34").Value = 4167.7144
$ws.Range("I34").Value = 4270.857
$ws.Range("J34").Value = 4098.952
$ws.Range("K34").Value = 4270.857
$ws.Range("L34").Value = 4098.952
$ws.Range("M34").Value = -4068.857
$ws.Range("N34").Value = -4502.952
$ws.Range("H50").Value = 29000
$ws.Range("J50").Value = 29000
$ws.Range("L50").Value = 29000
$ws.Range("N50").Value = -30250
$ws.Range("H51").Value = 28000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H59").Value = 40987.4
$ws.Range("J59").Value = 40987.4
$ws.Range("L59").Value = 40987.4
$ws.Range("N59").Value = -43277.4
$ws.Range("H60").Value = 27500
$ws.Range("J60").Value = 27500
$ws.Range("L60").Value = 27500
$ws.Range("N60").Value = -28522
$ws.Range("H61").Value = 28000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H68").Value = 41271.832
$ws.Range("J68").Value = 41271.832
$ws.Range("L68").Value = 41271.832
$ws.Range("N68").Value = -42769.832
$ws.Range("H71").Value = 41271.832
$ws.Range("J71").Value = 41271.832
$ws.Range("L71").Value = 123815.496
$ws.Range("N71").Value = -131303.496
$ws.Range("H107").Value = 487.48276
$ws.Range("I107").Value = 444
$ws.Range("K107").Value = 444
$ws.Range("M107").Value = 1476
$ws.Range("H122").Value = 4427.7856
$ws.Range("J122").Value = 4659.8
$ws.Range("L122").Value = 13979.4
$ws.Range("N122").Value = -18879.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9569.6875
$ws.Range("I3").Value = 8040.385
$ws.Range("K3").Value = 24121.155
$ws.Range("M3").Value = -24009.155
$ws.Range("H5").Value = 938.44446
$ws.Range("I5").Value = 491.5
$ws.Range("K5").Value = 1474.5
$ws.Range("M5").Value = -1362.5
$ws.Range("H135").Value = 938.44446
$ws.Range("I135").Value = 491.5
$ws.Range("K135").Value = 4423.5
$ws.Range("M135").Value = -1888.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 29776.111
$ws.Range("J46").Value = 38297
$ws.Range("L46").Value = 38297
$ws.Range("N46").Value = -38609
$ws.Range("H132").Value = 3167.2666
$ws.Range("J132").Value = 5557.4
$ws.Range("L132").Value = 16672.2
$ws.Range("N132").Value = -21732.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 21499.75
$ws.Range("I82").Value = 27666.334
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 27666.334
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -27305.334
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 21499.75
$ws.Range("I85").Value = 27666.334
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 27666.334
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -26418.334
$ws.Range("N85").Value = -5496
$ws.Range("H100").Value = 5277
$ws.Range("I100").Value = 5831
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 5831
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -5290
$ws.Range("N100").Value = -6082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 40999
$ws.Range("J51").Value = 40999
$ws.Range("L51").Value = 40999
$ws.Range("N51").Value = -42019
$ws.Range("H70").Value = 39576.668
$ws.Range("J70").Value = 40492
$ws.Range("L70").Value = 40492
$ws.Range("N70").Value = -41122
$ws.Range("H73").Value = 39576.668
$ws.Range("J73").Value = 40492
$ws.Range("L73").Value = 40492
$ws.Range("N73").Value = -42676
$ws.Range("H96").Value = 4214.4
$ws.Range("J96").Value = 6000
$ws.Range("L96").Value = 6000
$ws.Range("N96").Value = -8746
$ws.Range("H136").Value = 6393.8
$ws.Range("I136").Value = 5575.5835
$ws.Range("J136").Value = 9666.666999999999
$ws.Range("K136").Value = 16726.7505
$ws.Range("L136").Value = 29000.001
$ws.Range("M136").Value = -14176.7505
$ws.Range("N136").Value = -34100.001
